# Weekly fruit/vegetable price update:
# insert a new data row right after the header (new row 3), pushing the
# existing data rows down by one (old row 3 -> row 4, ..., old row 19 -> row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts rows 3..19 down to 4..20)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the latest weekly record
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 45037
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112013
$ws.Range("G3").Value = "Alcachofa"
$ws.Range("H3").Value = "Madrigal"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("N3").Value = "`$/caja 40 unidades"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 412
$ws.Range("Q3").Value = 40
$ws.Range("R3").Value = "Hortaliza"
